$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.505614041169197, 0.3375848360084654, 0.1529057820181812, 6.48142807727062, 8.477532736466465)
    3  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    4  = @(0.06328177979961902, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0.768386970581898)
    5  = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 6.48142807727062, 26.62400969366105)
    6  = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 10.35301142835362)
    7  = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 9.793184359356808)
    8  = @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.896700893398075)
    9  = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 8.260465185229014)
    10 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    11 = @(0.001754667048134761, 0.0001537489499301437, 3.082599426703578, 6.48142807727062, 9.565935919972263)
    12 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    13 = @(3.182878228561681, 1.65323645889881, 157.8057217802531, 0.4998867070740569, 163.1417231747877)
    14 = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 8.260465185229014)
    15 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    16 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    17 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    18 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
